$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1), columns D-H ---
# Columns D (old "Mean air temperature") is dropped; E/F, H/I, J/K, L/M pairs
# of "_x"/"_y" duplicate columns are collapsed into single renamed columns D-H.
$ws.Range("D1").Value = "4. Agriculture land area (% of land area)"
$ws.Range("E1").Value = "5. Average precipitation (mm per year)"
$ws.Range("F1").Value = "7. Fertilizer consumption (kilograms per hectare of arable land)"
$ws.Range("G1").Value = "13. Population"
$ws.Range("H1").Value = "17. Employment in agriculture (% of total employment) (modeled ILO estimate)"

# --- Remove now-unused trailing columns I:M (old duplicate "_y" columns) ---
$ws.Range("I1:M6").Delete()

# --- Replace data rows 2-8 (new data for years 2015-2021) ---
$ws.Cells.Item(2, 1).Value = "CPV"
$ws.Cells.Item(2, 2).Value = 2015
$ws.Cells.Item(2, 3).Value = 101.41
$ws.Cells.Item(2, 4).Value = 19.60297767
$ws.Cells.Item(2, 5).Value = 228
$ws.Cells.Item(2, 6).Value = 3.0346
$ws.Cells.Item(2, 7).Value = 552166
$ws.Cells.Item(2, 8).Value = 15.76179647527

$ws.Cells.Item(3, 1).Value = "CPV"
$ws.Cells.Item(3, 2).Value = 2016
$ws.Cells.Item(3, 3).Value = 98.06
$ws.Cells.Item(3, 4).Value = 19.60297767
$ws.Cells.Item(3, 5).Value = 228
$ws.Cells.Item(3, 6).Value = 3.0816
$ws.Cells.Item(3, 7).Value = 558394
$ws.Cells.Item(3, 8).Value = 14.7118964486198

$ws.Cells.Item(4, 1).Value = "CPV"
$ws.Cells.Item(4, 2).Value = 2017
$ws.Cells.Item(4, 3).Value = 91.26000000000001
$ws.Cells.Item(4, 4).Value = 19.60297767
$ws.Cells.Item(4, 5).Value = 228
$ws.Cells.Item(4, 6).Value = 3.9134
$ws.Cells.Item(4, 7).Value = 564954
$ws.Cells.Item(4, 8).Value = 13.6167403600925

$ws.Cells.Item(5, 1).Value = "CPV"
$ws.Cells.Item(5, 2).Value = 2018
$ws.Cells.Item(5, 3).Value = 72.04000000000001
$ws.Cells.Item(5, 4).Value = 19.60297767
$ws.Cells.Item(5, 5).Value = 228
$ws.Cells.Item(5, 6).Value = 4.07
$ws.Cells.Item(5, 7).Value = 571202
$ws.Cells.Item(5, 8).Value = 11.8382835657588

$ws.Cells.Item(6, 1).Value = "CPV"
$ws.Cells.Item(6, 2).Value = 2019
$ws.Cells.Item(6, 3).Value = 70.18000000000001
$ws.Cells.Item(6, 4).Value = 19.60297767
$ws.Cells.Item(6, 5).Value = 228
$ws.Cells.Item(6, 6).Value = 4.0876
$ws.Cells.Item(6, 7).Value = 577030
$ws.Cells.Item(6, 8).Value = 10.6044463385184

$ws.Cells.Item(7, 1).Value = "CPV"
$ws.Cells.Item(7, 2).Value = 2020
$ws.Cells.Item(7, 3).Value = 69.34
$ws.Cells.Item(7, 4).Value = 19.60297767
$ws.Cells.Item(7, 5).Value = 228
$ws.Cells.Item(7, 6).Value = 5.6728
$ws.Cells.Item(7, 7).Value = 582640
$ws.Cells.Item(7, 8).Value = 10.743090600802

$ws.Cells.Item(8, 1).Value = "CPV"
$ws.Cells.Item(8, 2).Value = 2021
$ws.Cells.Item(8, 3).Value = 70.52
$ws.Cells.Item(8, 4).Value = 19.60297767
$ws.Cells.Item(8, 5).Value = 228
$ws.Cells.Item(8, 6).Value = 5.373
$ws.Cells.Item(8, 7).Value = 587925
$ws.Cells.Item(8, 8).Value = 10.5378314239351

Write-Host "Edit complete"